$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignment (value does not parse as a plain number, so it
# is stored as text automatically).
function Set-PlainText($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

# Force text storage for values that otherwise look like plain numbers
# (e.g. "209.38") so Excel does not silently convert them to numeric
# cells / introduce floating point rounding.
function Set-ForcedText($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-PlainText 'D2' '28.289.18'
Set-PlainText 'E2' '  -0.90%  '
Set-PlainText 'D3' '1.552.63'
Set-PlainText 'E4' '  -0.02%  '
Set-ForcedText 'D5' '209.38'
Set-PlainText 'E5' '  -1.56%  '
Set-ForcedText 'D6' '0.484'
Set-PlainText 'E6' '  -1.61%  '
Set-PlainText 'E7' '  -0.03%  '
Set-ForcedText 'D8' '23.49'
Set-PlainText 'E8' '  -2.52%  '
Set-PlainText 'E9' '  -2.01%  '
Set-PlainText 'E10' '  -1.13%  '
Set-ForcedText 'D11' '0.0889'
Set-PlainText 'E11' '  +0.14%  '
Set-PlainText 'D12' '1.773.41'
Set-PlainText 'E12' '  -1.03%  '
Set-PlainText 'D13' '1.554.85'
Set-PlainText 'E13' '  -0.89%  '
Set-PlainText 'D14' '28.299.13'
Set-PlainText 'E14' '  -0.77%  '
Set-ForcedText 'D15' '3.63'
Set-PlainText 'E15' '  -1.24%  '
Set-ForcedText 'D16' '0.509'
Set-PlainText 'E16' '  -2.20%  '
Set-ForcedText 'D17' '60.48'
Set-PlainText 'E17' '  -2.77%  '
Set-ForcedText 'D18' '226.78'
Set-PlainText 'E18' '  -1.52%  '
Set-ForcedText 'D19' '7.31'
Set-PlainText 'E19' '  -0.52%  '
Set-PlainText 'D20' '0.0₃0675'
Set-PlainText 'E20' '  -2.37%  '
Set-PlainText 'E22' '  +1.25%  '
Set-ForcedText 'D23' '8.82'
Set-PlainText 'E23' '  -3.15%  '
Set-PlainText 'E24' '  -5.49%  '
Set-ForcedText 'D25' '147.87'
Set-PlainText 'E25' '  -2.20%  '
Set-ForcedText 'D26' '14.76'
Set-PlainText 'E26' '  -1.70%  '
Set-PlainText 'E27' '  -0.22%  '
Set-PlainText 'E28' '  -0.06%  '
Set-ForcedText 'D29' '6.22'
Set-PlainText 'E29' '  -3.03%  '
Set-PlainText 'E30' '  -3.51%  '
Set-PlainText 'E31' '  -4.24%  '
Set-PlainText 'E32' '  -0.87%  '
Set-PlainText 'E33' '  -1.06%  '
Set-PlainText 'D34' '1.385.16'
Set-PlainText 'E34' '  -0.53%  '
Set-PlainText 'E35' '  +0.59%  '
Set-PlainText 'E36' '  -2.60%  '
Set-PlainText 'E37' '  -1.37%  '
Set-ForcedText 'D38' '2.58'
Set-PlainText 'E38' '  -1.45%  '
Set-PlainText 'E39' '  -2.20%  '
Set-ForcedText 'D40' '1.94'
Set-PlainText 'E40' '  +2.40%  '
Set-ForcedText 'D41' '0.512'
Set-PlainText 'E41' '  -2.05%  '
Set-PlainText 'E42' '  -0.07%  '
Set-PlainText 'E43' '  -1.39%  '
Set-ForcedText 'D44' '0.0467'
Set-PlainText 'E44' '  +1.15%  '
Set-PlainText 'E45' '  -1.19%  '
Set-ForcedText 'D46' '61.79'
Set-PlainText 'E46' '  -1.61%  '
Set-PlainText 'D47' '1.687.33'
Set-PlainText 'E47' '  -1.02%  '
Set-PlainText 'E48' '  -6.68%  '
Set-ForcedText 'D49' '85.53'
Set-PlainText 'E49' '  -0.94%  '
Set-ForcedText 'D50' '42.18'
Set-PlainText 'E50' '  +6.53%  '
Set-PlainText 'E51' '  +0.29%  '
